$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.874.40"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.563.76"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0583"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.787.52"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.564.76"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "26.887.42"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "1.403.25"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.918"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "1.700.29"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0504"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "0.0₇0981"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0945"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
